$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 1254
$ws.Range("I2").Value = 3351
$ws.Range("J2").Value = 14101
$ws.Range("K2").Value = 65
$ws.Range("L2").Value = 3784
$ws.Range("M2").Value = 240
$ws.Range("N2").Value = 2521
$ws.Range("O2").Value = 6
$ws.Range("P2").Value = 49
$ws.Range("Q2").Value = 18
$ws.Range("R2").Value = 193
$ws.Range("S2").Value = 1582
$ws.Range("T2").Value = 2502
$ws.Range("U2").Value = 202
$ws.Range("V2").Value = 21929
$ws.Range("X2").Value = 22109
$ws.Range("Y2").Value = 44
$ws.Range("Z2").Value = 327
$ws.Range("AA2").Value = 130
